## Generate Report for Handback
## Populates the "Latest Target File" / "Latest Handback File" /
## "Latest Handback DateTime" / "Error Detail" columns (I, J, K, P) for the
## last row of the zh-cn and de-de sheets, widens the Error Detail column,
## and links the new "Latest Target File" cell to the handoff markdown file
## (mirroring the existing link in column A).

$wb = $excel.ActiveWorkbook

# -- zh-cn sheet (row 8) -----------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P) to match the other wide columns.
$wsZh.Columns.Item(16).ColumnWidth = $wsZh.Columns.Item(1).ColumnWidth

# Find the handoff-file hyperlink target already used in column A so the
# new "Latest Target File" link points at the very same commit.
$zhTarget = $null
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$8") {
        $zhTarget = $hl.Address
    }
}

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $zhTarget, [Type]::Missing, [Type]::Missing, "cc82e517-e93e-4749-b293-117f93885f13.md")

$wsZh.Range("J8").Value = "cc82e517-e93e-4749-b293-117f93885f13.a65aa8a853960c7482d302306efe1b01400cea9d.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-21 20:54:41"
$wsZh.Range("P8").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/470e5796d3087b457a8d2c0f385e07cb9f11bcbe/e2e/cc82e517-e93e-4749-b293-117f93885f13.md, latest: $zhTarget."

# -- de-de sheet (row 8) -------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = $wsDe.Columns.Item(1).ColumnWidth

$deTarget = $null
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$8") {
        $deTarget = $hl.Address
    }
}

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $deTarget, [Type]::Missing, [Type]::Missing, "cc82e517-e93e-4749-b293-117f93885f13.md")

$wsDe.Range("J8").Value = "cc82e517-e93e-4749-b293-117f93885f13.a65aa8a853960c7482d302306efe1b01400cea9d.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-21 20:54:48"
$wsDe.Range("P8").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/470e5796d3087b457a8d2c0f385e07cb9f11bcbe/e2e/cc82e517-e93e-4749-b293-117f93885f13.md, latest: $deTarget."
